$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row: "_old" -> "_FV2404" and "_new" -> "_FV2410"
#    (the "diff" header in column K stays unchanged)
$basenames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $basenames.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $basenames[$i] + "_FV2404"
}

for ($i = 0; $i -lt $basenames.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $basenames[$i] + "_FV2410"
}

# 2. Turn the populated range into an Excel Table ("Table1") with autofilter
$range = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (split below row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Edit complete"
